$d = $word.ActiveDocument

# The last paragraph in the document body (right before the sectPr) is
# an empty paragraph. The edit appends new italic text to it, split
# across several runs:
#   "TBD " / "how to " / "add" / "/" /
#   "append to existing commit (saw it on a YT vid earlier)"

$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)

# Mark the (currently empty) paragraph mark itself as italic first, so
# the paragraph's own mark formatting (w:pPr/w:rPr) picks up
# <w:i/><w:iCs/> exactly like the target markup.
$mark = $target.Range
$mark.Font.Italic = $true
$mark.Font.ItalicBi = $true

$chunks = @(
    "TBD ",
    "how to ",
    "add",
    "/",
    "append to existing commit (saw it on a YT vid earlier)"
)

foreach ($chunk in $chunks) {
    $r = $target.Range
    $r.InsertAfter($chunk)
    $r.Font.Italic = $true
    $r.Font.ItalicBi = $true
    $r.LanguageID = "en-US"
}

Write-Output "Appended TBD note to final paragraph."
